$wb = $excel.ActiveWorkbook

# Add the missing "Opex" column to the Sinks sheet (between Capex and
# EnergyType), matching the Sources and Transformers sheets which already
# have this column.
$wsSinks = $wb.Worksheets.Item("Sinks")
$wsSinks.Range("C1").EntireColumn.Insert()
$wsSinks.Range("C1").Value = "Opex"

# Move the active tab from Sources to Connectors.
$wsConnectors = $wb.Worksheets.Item("Connectors")
$wsConnectors.Activate()
